$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up row 462 (last existing row) ---
# Remove the leftover empty placeholder cells that used to come from the
# Google Form submission (C, D, G, I were blank inline strings).
$ws.Range("C462").ClearContents()
$ws.Range("D462").ClearContents()
$ws.Range("G462").ClearContents()
$ws.Range("I462").ClearContents()

# Re-stamp the submission timestamp for row 462 with the refreshed value.
$ws.Range("L462").Value = 45281.72770226852

# --- Make sure the new timestamp cells use the same date/time format ---
# as the rest of column L (numFmt "yyyy-mm-dd h:mm:ss").
$ws.Range("L463:L476").NumberFormat = "yyyy-mm-dd h:mm:ss"

# --- New penalty log entries pulled in via the Google Forms API ---

# Row 463
$ws.Range("A463").Value = "'6359"
$ws.Range("B463").Value = "_gorlomi"
$ws.Range("E463").Value = "CL"
$ws.Range("F463").Value = 3
$ws.Range("H463").Value = "https://youtu.be/&t=3525"
$ws.Range("J463").Value = "_gorlomi"
$ws.Range("K463").Value = "edfvswfv"
$ws.Range("L463").Value = 45310.00884177083

# Row 464
$ws.Range("A464").Value = "'6359"
$ws.Range("B464").Value = "_gorlomi"
$ws.Range("C464").Value = "fgnhfgnh"
$ws.Range("D464").Value = "dhbfth"
$ws.Range("E464").Value = "CL"
$ws.Range("F464").Value = 3
$ws.Range("G464").Value = "Edging/Bumping/Swerving on Straights"
$ws.Range("H464").Value = "https://youtu.be/&t=3525"
$ws.Range("I464").Value = "dhttdrh"
$ws.Range("J464").Value = "_gorlomi"
$ws.Range("K464").Value = "edfvswfv"
$ws.Range("L464").Value = 45310.00884177083

# Row 465
$ws.Range("A465").Value = "'9884"
$ws.Range("B465").Value = "_gorlomi"
$ws.Range("E465").Value = "UL"
$ws.Range("F465").Value = 3
$ws.Range("H465").Value = "https://youtu.be/&t=234"
$ws.Range("J465").Value = "_gorlomi"
$ws.Range("K465").Value = "eafvf"
$ws.Range("L465").Value = 45310.0130899537

# Row 466
$ws.Range("A466").Value = "'9884"
$ws.Range("B466").Value = "_gorlomi"
$ws.Range("C466").Value = "No offence"
$ws.Range("E466").Value = "UL"
$ws.Range("F466").Value = 3
$ws.Range("H466").Value = "https://youtu.be/&t=234"
$ws.Range("I466").Value = "rdgdr"
$ws.Range("J466").Value = "_gorlomi"
$ws.Range("K466").Value = "eafvf"
$ws.Range("L466").Value = 45310.0130899537

# Row 467
$ws.Range("A467").Value = "'6437"
$ws.Range("B467").Value = "_gorlomi"
$ws.Range("E467").Value = "UL"
$ws.Range("F467").Value = 5
$ws.Range("H467").Value = "https://youtu.be/&t=352"
$ws.Range("J467").Value = "_gorlomi"
$ws.Range("K467").Value = "gvsrgv"
$ws.Range("L467").Value = 45310.01524820602

# Row 468
$ws.Range("A468").Value = "'6437"
$ws.Range("B468").Value = "_gorlomi"
$ws.Range("C468").Value = "No offence"
$ws.Range("E468").Value = "UL"
$ws.Range("F468").Value = 5
$ws.Range("H468").Value = "https://youtu.be/&t=352"
$ws.Range("I468").Value = "rhbtrb"
$ws.Range("J468").Value = "_gorlomi"
$ws.Range("K468").Value = "gvsrgv"
$ws.Range("L468").Value = 45310.01524820602

# Row 469
$ws.Range("A469").Value = "'5446"
$ws.Range("B469").Value = "_gorlomi"
$ws.Range("E469").Value = "UL"
$ws.Range("F469").Value = 3
$ws.Range("H469").Value = "https://youtu.be/&t=352"
$ws.Range("J469").Value = "_gorlomi"
$ws.Range("K469").Value = "edgfw"
$ws.Range("L469").Value = 45310.01578501157

# Row 470
$ws.Range("A470").Value = "'5446"
$ws.Range("B470").Value = "_gorlomi"
$ws.Range("C470").Value = "No offence"
$ws.Range("E470").Value = "UL"
$ws.Range("F470").Value = 3
$ws.Range("H470").Value = "https://youtu.be/&t=352"
$ws.Range("I470").Value = "fvdsvs"
$ws.Range("J470").Value = "_gorlomi"
$ws.Range("K470").Value = "edgfw"
$ws.Range("L470").Value = 45310.01578501157

# Row 471
$ws.Range("A471").Value = "'3925"
$ws.Range("B471").Value = "_gorlomi"
$ws.Range("E471").Value = "UL"
$ws.Range("F471").Value = 4
$ws.Range("H471").Value = "https://youtu.be/&t=352"
$ws.Range("J471").Value = "_gorlomi"
$ws.Range("K471").Value = "dafvadsv"
$ws.Range("L471").Value = 45310.01768652778

# Row 472
$ws.Range("A472").Value = "'3925"
$ws.Range("B472").Value = "_gorlomi"
$ws.Range("C472").Value = "No offence"
$ws.Range("E472").Value = "UL"
$ws.Range("F472").Value = 4
$ws.Range("H472").Value = "https://youtu.be/&t=352"
$ws.Range("I472").Value = "afcsdv"
$ws.Range("J472").Value = "_gorlomi"
$ws.Range("K472").Value = "dafvadsv"
$ws.Range("L472").Value = 45310.01768652778

# Row 473
$ws.Range("A473").Value = "'6747"
$ws.Range("B473").Value = "_gorlomi"
$ws.Range("E473").Value = "UL"
$ws.Range("F473").Value = 4
$ws.Range("H473").Value = "https://youtu.be/&t=352"
$ws.Range("J473").Value = "_gorlomi"
$ws.Range("K473").Value = "svdrg"
$ws.Range("L473").Value = 45310.02423076389

# Row 474
$ws.Range("A474").Value = "'6747"
$ws.Range("B474").Value = "_gorlomi"
$ws.Range("C474").Value = "No offence"
$ws.Range("E474").Value = "UL"
$ws.Range("F474").Value = 4
$ws.Range("H474").Value = "https://youtu.be/&t=352"
$ws.Range("I474").Value = "fgdrg"
$ws.Range("J474").Value = "_gorlomi"
$ws.Range("K474").Value = "svdrg"
$ws.Range("L474").Value = 45310.02423076389

# Row 475
$ws.Range("A475").Value = "'2171"
$ws.Range("B475").Value = "_gorlomi"
$ws.Range("E475").Value = "CL"
$ws.Range("F475").Value = 4
$ws.Range("H475").Value = "https://youtu.be/&t=352"
$ws.Range("J475").Value = "_gorlomi"
$ws.Range("K475").Value = "wesgvswgv"
$ws.Range("L475").Value = 45310.02582716435

# Row 476
$ws.Range("A476").Value = "'2171"
$ws.Range("B476").Value = "_gorlomi"
$ws.Range("C476").Value = "No offence"
# D476 stays blank (empty form answer) but the cell itself still gets
# written out, so just touch it without putting any style on it.
$ws.Range("D476").Style = "Normal"
$ws.Range("E476").Value = "CL"
$ws.Range("F476").Value = 4
# Same idea for G476 - a blank answer that still materializes the cell.
$ws.Range("G476").Style = "Normal"
$ws.Range("H476").Value = "https://youtu.be/&t=352"
$ws.Range("I476").Value = "dasvfds"
$ws.Range("J476").Value = "_gorlomi"
$ws.Range("K476").Value = "wesgvswgv"
$ws.Range("L476").Value = 45310.02582716673
